$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mirror column D's formatting into the new column E (rows 3-8) by
# copying the existing formats, then fill in the new data values.
$ws.Range("D3:D8").Copy()
$ws.Range("E3:E8").PasteSpecial(-4122)

$ws.Range("E4").Value = 2020
$ws.Range("E5").Value = 11.5
$ws.Range("E6").Value = 2.6
$ws.Range("E7").Value = 2
$ws.Range("E8").Value = 0.3

# E7 gets its own one-decimal number format (new style), matching the
# rest of the percentage figures in the table.
$ws.Range("E7").NumberFormat = "0.0"

# Restore the cursor/selection that was active when the file was saved.
$ws.Range("B15").Select()
